# Bidirectional Excel sync PoC: update the BLOCKCHAINS row for avalanche
# with a newly-synced test row (new id, refreshed RPC endpoints, gas/fee
# config, contract addresses, dex/protocol lists, and live health stats).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BLOCKCHAINS")

# Identity
$ws.Range("A5").Value = "avalanche_test_1760873174"

# RPC / network endpoints
$ws.Range("F5").Value = "https://avalanche.public-rpc.com"
$ws.Range("H5").Value = "https://rpc.ankr.com/avalanche"
$ws.Range("I5").Value = "wss://avalanche.public-rpc.com"

# Gas / fee configuration
$ws.Range("K5").Value = 2000
$ws.Range("L5").Value = 25
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 25
$ws.Range("O5").Value = $true
$ws.Range("P5").Value = 25
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 15000000

# Contract addresses
$ws.Range("S5").Value = "0xcA11bde05977b3631167028862bE2a173976CA11"
$ws.Range("T5").Value = "0x49D5c2BdFfac6CE2BFdB6640F4F80f226bc10bAB"
$ws.Range("U5").Value = "0xB97EF9Ef8734C71904D8002F8b6Bc66Dd9c48a6E"
$ws.Range("V5").Value = "0x9702230A8Ea53601f5cD2dc00fDBc13d4dF4A8c7"
$ws.Range("W5").Value = "0xd586E7F844cEa2F87f50152665BCbc2C279D8d70"

# Supported dexes / protocols
$ws.Range("X5").Value = "TraderJoe, Pangolin, SushiSwap, Curve"
$ws.Range("Y5").Value = "Aave, Benqi, TraderJoe"

# Gas/finality/protection flags
$ws.Range("AC5").Value = 21000
$ws.Range("AD5").Value = 1
$ws.Range("AE5").Value = $true
$ws.Range("AF5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AH5").Value = $false

# Live health / sync status
$ws.Range("AI5").Value = "degraded"
$ws.Range("AL5").Value = "synced"
$ws.Range("AN5").Value = 99.5
$ws.Range("AO5").Value = 0.5
$ws.Range("AP5").Value = 0
$ws.Range("AQ5").Value = 15000
$ws.Range("AR5").Value = 3
$ws.Range("AS5").Value = 10
$ws.Range("AT5").Value = 100
